# Test/Test.xlsx was re-uploaded with a refreshed 5-row (time, input) data
# table: the old bold "time"/"input" header row is gone, and column A/B now
# hold a plain numeric time-step / input-index series starting at row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the old header + 4 data rows with the new 5-row numeric table.
# Column A: 0.1, 1.1, 2.1, 3.1, 4.1   Column B: 1, 2, 3, 4, 5
for ($i = 0; $i -lt 5; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = 0.1 + $i
    $ws.Cells.Item($r, 2).Value = $i + 1
}

# The former header row (A1:B1) was bold text ("time"/"input"); now that the
# cells hold plain numbers, drop the bold formatting too.
$ws.Range("A1:B1").ClearFormats()

# Match the saved selection/active cell (just below the new table).
$ws.Range("A6").Select() | Out-Null
